$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Good Morning" greeting text with "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update the selection/active cell shown in the sheet view to E8
$ws.Range("E8").Select()
